# Inserts a new data row at row 146 (pushing existing rows 146-244 down to
# 147-245) and populates it with the new observation. The rest of the
# columns (A,B,C,E,F,G,H,I,J,N,Q,R) are identical to the row that used to be
# at 146 (now 147), so after the insert they already carry over correctly;
# we just need to set the handful of cells that actually hold new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 146:244 down to 147:245 by inserting a new row at 146.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new record.
$ws.Cells.Item(146, 1).Value = 9
$ws.Cells.Item(146, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(146, 3).Value = "Metropolitana"
$ws.Cells.Item(146, 4).Value = 44767
$ws.Cells.Item(146, 5).Value = 13
$ws.Cells.Item(146, 6).Value = 100112026
$ws.Cells.Item(146, 7).Value = "Haba"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 52
$ws.Cells.Item(146, 11).Value = 18000
$ws.Cells.Item(146, 12).Value = 18000
$ws.Cells.Item(146, 13).Value = 18000
$ws.Cells.Item(146, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(146, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(146, 16).Value = 720
$ws.Cells.Item(146, 17).Value = 25
$ws.Cells.Item(146, 18).Value = "Hortaliza"
